$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Équipe Maxime Cliche"
$ws.Range("B4").Select()
